# Updated with Excel support
# Applies the "virtualan_collection_pet" workbook edit:
#  - PetPost (row 2) test case is repointed at a mockbin.org endpoint, switches
#    from JSON to XML payloads/files, gains an ExcludeField value, and its
#    expected HTTP status moves from 201 to 200.
#  - A new ExcludeField header value is written to H2 ("Date").
#  - Column B is widened slightly.
#  - The PetPost URL cell (still hyperlinked to the original localhost target)
#    gets a darker link-text color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (PetPost) content updates -------------------------------------

# URL text is replaced; the existing hyperlink (still targeting
# http://localhost:8800/api/pets) is preserved because only the cell's
# displayed text is being changed, not the hyperlink target itself.
$ws.Range("C2").Value = "http://mockbin.org/bin/2c5f64fe-4b65-4453-85a5-5308767e79e8"

$ws.Range("D2").Value = "application/xml"
$ws.Range("E2").Value = "input.xml"
$ws.Range("G2").Value = "output.xml"

# New ExcludeField cell for the PetPost row.
$ws.Range("H2").Value = "VirtualanStdType=EDI-271"

$ws.Range("J2").Value = "Date"
$ws.Range("K2").Value = 200

# --- Column sizing ----------------------------------------------------------
$ws.Range("B:B").ColumnWidth = 17.29

# --- Visual tweak: darken the PetPost hyperlink cell's font color ----------
$ws.Range("C2").Font.Color = 5263440
